$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (A1:E10) was moved down one row and right one column, to B2:F11.
# Insert a blank row above row 1 and a blank column to the left of column A
# so the existing table content/styles land on B2:F11.
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# Row-height bookkeeping: row heights are a row-level (not cell-level) property,
# so after the move the two distinct heights used by the table need to be
# re-applied at their final row positions.
$ws.Rows.Item(3).RowHeight = 13.5
$ws.Rows.Item(4).RowHeight = 15.75

# Column-width bookkeeping: likewise a column-level property. Column A keeps its
# original width, and columns B:F pick up new (best-fit) widths for the moved data.
$ws.Columns.Item(1).ColumnWidth = 33.42578125
$ws.Columns.Item(2).ColumnWidth = 34.7109375
$ws.Columns.Item(3).ColumnWidth = 9.140625
$ws.Columns.Item(4).ColumnWidth = 11.85546875
$ws.Columns.Item(5).ColumnWidth = 10.28515625
$ws.Columns.Item(6).ColumnWidth = 12.7109375

# Apply a thin border around every cell of the table (B2:F11).
$ws.Range("B2:F11").Borders.LineStyle = 1
$ws.Range("B2:F11").Borders.Weight = 2

# Match the author's final selection.
$ws.Range("H10").Select()
